$d = $word.ActiveDocument

# "Versi" + "on"  ->  merge into a single "Version" run
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2)

# " 2" -> " 1."  (the run right before the _GoBack bookmark)
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 1.", 2)

# Drop the now-redundant trailing "." run that used to sit after the bookmark
$bm = $d.Bookmarks("_GoBack")
$tail = $d.Range($bm.End, $d.Content.End - 1)
$tail.Delete()
